$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their original text formatting so that
# numeric-looking strings (e.g. "318.67") are not silently converted to
# floating point numbers by Excel's automatic type inference.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '42.793.36'
$ws.Range('E2').Value = '  +1.43%  '
$ws.Range('D3').Value = '2.317.01'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '318.67'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('D6').Value = '105.32'
$ws.Range('E6').Value = '  +1.52%  '
$ws.Range('D7').Value = '0.631'
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.612'
$ws.Range('E9').Value = '  +0.96%  '
$ws.Range('D10').Value = '40.26'
$ws.Range('E10').Value = '  +2.01%  '
$ws.Range('D11').Value = '0.0911'
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').Value = '8.55'
$ws.Range('E12').Value = '  +3.38%  '
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('D14').Value = '0.989'
$ws.Range('E14').Value = '  +2.74%  '
$ws.Range('D15').Value = '15.54'
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('D16').Value = '2.670.06'
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('D17').Value = '2.317.58'
$ws.Range('E17').Value = '  +1.12%  '
$ws.Range('D18').Value = '42.624.24'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').Value = '7.68'
$ws.Range('E19').Value = '  +2.93%  '
$ws.Range('D20').Value = '0.0000107'
$ws.Range('E20').Value = '  +0.89%  '
$ws.Range('D21').Value = '13.52'
$ws.Range('E21').Value = '  +35.78%  '
$ws.Range('D22').Value = '74.07'
$ws.Range('E22').Value = '  +0.95%  '
$ws.Range('E23').Value = '  -1.22%  '
$ws.Range('D24').Value = '270.27'
$ws.Range('E24').Value = '  -3.12%  '
$ws.Range('D25').Value = '2.25'
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('D27').Value = '10.96'
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('D28').Value = '2.28'
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('D29').Value = '22.76'
$ws.Range('E29').Value = '  -0.76%  '
$ws.Range('D30').Value = '38.49'
$ws.Range('E30').Value = '  +5.21%  '
$ws.Range('D31').Value = '6.49'
$ws.Range('E31').Value = '  +11.40%  '
$ws.Range('D32').Value = '166.91'
$ws.Range('E32').Value = '  +2.30%  '
$ws.Range('D33').Value = '0.0891'
$ws.Range('E33').Value = '  +1.87%  '
$ws.Range('D34').Value = '0.133'
$ws.Range('E34').Value = '  -2.54%  '
$ws.Range('D35').Value = '2.62'
$ws.Range('E35').Value = '  -8.03%  '
$ws.Range('E36').Value = '  +0.86%  '
$ws.Range('D37').Value = '4.61'
$ws.Range('E37').Value = '  +1.34%  '
$ws.Range('D38').Value = '0.0355'
$ws.Range('E38').Value = '  +1.28%  '
$ws.Range('D39').Value = '3.74'
$ws.Range('E39').Value = '  +0.76%  '
$ws.Range('D40').Value = '2.81'
$ws.Range('E40').Value = '  -4.41%  '
$ws.Range('D41').Value = '1.66'
$ws.Range('E41').Value = '  +14.07%  '
$ws.Range('D42').Value = '101.40'
$ws.Range('E42').Value = '  +2.21%  '
$ws.Range('D43').Value = '70.89'
$ws.Range('E43').Value = '  +2.05%  '
$ws.Range('D44').Value = '0.227'
$ws.Range('E44').Value = '  +1.07%  '
$ws.Range('E45').Value = '  -0.47%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '117.91'
$ws.Range('E46').Value = '  +3.95%  '
$ws.Range('B47').Value = 'Celestia'
$ws.Range('C47').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D47').Value = '12.42'
$ws.Range('E47').Value = '  +4.03%  '
$ws.Range('D48').Value = '82.13'
$ws.Range('E48').Value = '  +6.59%  '
$ws.Range('D49').Value = '1.648.37'
$ws.Range('E49').Value = '  +4.05%  '
$ws.Range('D50').Value = '5.34'
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('D51').Value = '8.89'
$ws.Range('E51').Value = '  -1.21%  '
